# Update cryptocurrency price/volume data per upstream refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.873.90"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "1.642.41"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").Value = "'215.85"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").Value = "'0.5061"
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("D7").Value = "'1.005"
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("D8").Value = "'0.2584"
$ws.Range("E8").Value = "  +0.19%  "
$ws.Range("D9").Value = "'0.06434"
$ws.Range("E9").Value = "  +1.50%  "
$ws.Range("D10").Value = "'20.44"
$ws.Range("E10").Value = "  +4.96%  "
$ws.Range("D11").Value = "'0.07808"
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("D12").Value = "'4.268"
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "1.871.70"
$ws.Range("E13").Value = "  +0.65%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.640.63"
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("D15").Value = "'0.5626"
$ws.Range("E15").Value = "  +2.14%  "
$ws.Range("D16").Value = "0.0₅7689"
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").Value = "'63.30"
$ws.Range("E17").Value = "  -0.95%  "
$ws.Range("D18").Value = "25.911.62"
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("E19").Value = "  +0.43%  "
$ws.Range("D20").Value = "'192.81"
$ws.Range("E20").Value = "  -1.20%  "
$ws.Range("D21").Value = "'4.370"
$ws.Range("E21").Value = "  -1.01%  "
$ws.Range("D22").Value = "'9.926"
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("D23").Value = "'6.116"
$ws.Range("E23").Value = "  +1.09%  "
$ws.Range("D24").Value = "'1.005"
$ws.Range("E24").Value = "  +0.30%  "
$ws.Range("D25").Value = "'1.802"
$ws.Range("E25").Value = "  -6.11%  "
$ws.Range("D26").Value = "'141.20"
$ws.Range("E26").Value = "  -0.98%  "
$ws.Range("D27").Value = "'0.1238"
$ws.Range("E27").Value = "  -1.10%  "
$ws.Range("D28").Value = "'6.788"
$ws.Range("E28").Value = "  +0.28%  "
$ws.Range("D29").Value = "'15.53"
$ws.Range("D30").Value = "'1.245"
$ws.Range("E30").Value = "  +0.44%  "
$ws.Range("D31").Value = "'0.04934"
$ws.Range("E31").Value = "  +0.61%  "
$ws.Range("D32").Value = "'3.292"
$ws.Range("E32").Value = "  +1.29%  "
$ws.Range("D33").Value = "'3.237"
$ws.Range("E33").Value = "  +1.26%  "
$ws.Range("D34").Value = "'1.575"
$ws.Range("E34").Value = "  +2.12%  "
$ws.Range("D35").Value = "'2.389"
$ws.Range("E35").Value = "  +0.84%  "
$ws.Range("D36").Value = "'0.9052"
$ws.Range("E36").Value = "  +0.71%  "
$ws.Range("D37").Value = "'0.5545"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("D38").Value = "1.132.49"
$ws.Range("E38").Value = "  +1.44%  "
$ws.Range("D39").Value = "'2.554"
$ws.Range("E39").Value = "  +0.60%  "
$ws.Range("E40").Value = "  +0.65%  "
$ws.Range("D41").Value = "'1.003"
$ws.Range("E41").Value = "  +0.28%  "
$ws.Range("D42").Value = "'5.515"
$ws.Range("E42").Value = "  -1.86%  "
$ws.Range("D43").Value = "'0.8029"
$ws.Range("E43").Value = "  +0.74%  "
$ws.Range("D44").Value = "'98.95"
$ws.Range("D45").Value = "1.782.29"
$ws.Range("E45").Value = "  +0.63%  "
$ws.Range("E46").Value = "  -5.99%  "
$ws.Range("D47").Value = "'55.69"
$ws.Range("E47").Value = "  +1.64%  "
$ws.Range("D48").Value = "'0.4292"
$ws.Range("E48").Value = "  -3.46%  "
$ws.Range("D49").Value = "'7.760"
$ws.Range("E49").Value = "  +2.42%  "
$ws.Range("D50").Value = "'0.05046"
$ws.Range("E50").Value = "  -1.71%  "
$ws.Range("D51").Value = "'0.9989"
$ws.Range("E51").Value = "  -0.55%  "
